$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.649.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.598.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.592.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.629.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("E23").Value = "  -3.51%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.666"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.295.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.734.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.892"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("E49").Value = "  +2.12%  "

$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("E51").Value = "  +0.39%  "
